$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 276.125
$ws.Range("I12").Value = 210
$ws.Range("J12").Value = 386.33334
$ws.Range("K12").Value = 210
$ws.Range("L12").Value = 386.33334
$ws.Range("M12").Value = -40
$ws.Range("N12").Value = -726.33334
$ws.Range("H17").Value = 2712.375
$ws.Range("J17").Value = 2416.5
$ws.Range("L17").Value = 7249.5
$ws.Range("N17").Value = -7585.5
$ws.Range("H32").Value = 27778478
$ws.Range("I32").Value = 41667170
$ws.Range("K32").Value = 41667170
$ws.Range("M32").Value = -41666844
$ws.Range("H40").Value = 3727.111
$ws.Range("I40").Value = 3934.8572
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 3934.8572
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -3759.8572
$ws.Range("N40").Value = -3350
$ws.Range("H87").Value = 114996.336
$ws.Range("J87").Value = 114996.336
$ws.Range("L87").Value = 114996.336
$ws.Range("N87").Value = -117492.336
$ws.Range("H90").Value = 114996.336
$ws.Range("J90").Value = 114996.336
$ws.Range("L90").Value = 344989.008
$ws.Range("N90").Value = -357469.008
$ws.Range("H112").Value = 2187.2856
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2187.2856
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 6561.8568
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -8777.856800000001
$ws.Range("H132").Value = 7543.409
$ws.Range("I132").Value = 7797.8
$ws.Range("K132").Value = 23393.4
$ws.Range("M132").Value = -20863.4
$ws.Range("H137").Value = 1222760.9
$ws.Range("I137").Value = 6252137.5
$ws.Range("K137").Value = 18756412.5
$ws.Range("M137").Value = -18753862.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 507856.72
$ws.Range("I74").Value = 1852784.6
$ws.Range("J74").Value = 3508.75
$ws.Range("K74").Value = 1852784.6
$ws.Range("L74").Value = 3508.75
$ws.Range("M74").Value = -1851910.6
$ws.Range("N74").Value = -5256.75
$ws.Range("H77").Value = 507856.72
$ws.Range("I77").Value = 1852784.6
$ws.Range("J77").Value = 3508.75
$ws.Range("K77").Value = 9263923
$ws.Range("L77").Value = 17543.75
$ws.Range("M77").Value = -9259555
$ws.Range("N77").Value = -26279.75
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H132").Value = 5785.2856
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5785.2856
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 17355.8568
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -22415.8568
$ws.Range("H133").Value = 102995
$ws.Range("J133").Value = 102995
$ws.Range("L133").Value = 102995
$ws.Range("N133").Value = -108055

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3112.7073
$ws.Range("I134").Value = 3000.8276
$ws.Range("J134").Value = 3383.0833
$ws.Range("K134").Value = 9002.4828
$ws.Range("L134").Value = 10149.2499
$ws.Range("M134").Value = -6467.4828
$ws.Range("N134").Value = -15219.2499
$ws.Range("H139").Value = 23923.924
$ws.Range("J139").Value = 23923.924
$ws.Range("L139").Value = 23923.924
$ws.Range("N139").Value = -34203.924

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 95788.5
$ws.Range("J52").Value = 95788.5
$ws.Range("L52").Value = 95788.5
$ws.Range("N52").Value = -96376.5
$ws.Range("H99").Value = 4999.5
$ws.Range("I99").Value = 4999.5
$ws.Range("K99").Value = 4999.5
$ws.Range("M99").Value = -3501.5
$ws.Range("H122").Value = 4623.5625
$ws.Range("I122").Value = 4183.1816
$ws.Range("J122").Value = 4854.2383
$ws.Range("K122").Value = 12549.5448
$ws.Range("L122").Value = 14562.7149
$ws.Range("M122").Value = -10099.5448
$ws.Range("N122").Value = -19462.7149
$ws.Range("H126").Value = 4999.5
$ws.Range("I126").Value = 4999.5
$ws.Range("K126").Value = 14998.5
$ws.Range("M126").Value = -12528.5
$ws.Range("H134").Value = 2270.5625
$ws.Range("I134").Value = 2022.6
$ws.Range("J134").Value = 3156.1428
$ws.Range("K134").Value = 6067.799999999999
$ws.Range("L134").Value = 9468.428400000001
$ws.Range("M134").Value = -3532.799999999999
$ws.Range("N134").Value = -14538.4284
$ws.Range("H138").Value = 68887
$ws.Range("J138").Value = 68887
$ws.Range("L138").Value = 68887
$ws.Range("N138").Value = -79167

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1026.3334
$ws.Range("J5").Value = 1289.375
$ws.Range("L5").Value = 3868.125
$ws.Range("N5").Value = -4092.125
$ws.Range("H36").Value = 1466.3334
$ws.Range("I36").Value = 574.5
$ws.Range("J36").Value = 3250
$ws.Range("K36").Value = 1723.5
$ws.Range("L36").Value = 9750
$ws.Range("M36").Value = -1554.5
$ws.Range("N36").Value = -10088
$ws.Range("H121").Value = 233798.89
$ws.Range("I121").Value = 200
$ws.Range("J121").Value = 300541.44
$ws.Range("K121").Value = 600
$ws.Range("L121").Value = 901624.3200000001
$ws.Range("M121").Value = 710
$ws.Range("N121").Value = -904244.3200000001
$ws.Range("H135").Value = 1026.3334
$ws.Range("J135").Value = 1289.375
$ws.Range("L135").Value = 11604.375
$ws.Range("N135").Value = -16674.375
$ws.Range("H140").Value = 14208.821
$ws.Range("I140").Value = 11399.6
$ws.Range("K140").Value = 34198.8
$ws.Range("M140").Value = -29018.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3410.72
$ws.Range("I132").Value = 2888.5
$ws.Range("J132").Value = 5499.6
$ws.Range("K132").Value = 8665.5
$ws.Range("L132").Value = 16498.8
$ws.Range("M132").Value = -6135.5
$ws.Range("N132").Value = -21558.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 7007.0835
$ws.Range("J7").Value = 9499.5
$ws.Range("K7").Value = 7007.0835
$ws.Range("L7").Value = 9499.5
$ws.Range("M7").Value = -6895.0835
$ws.Range("N7").Value = -9723.5
$ws.Range("H40").Value = 4733.8057
$ws.Range("I40").Value = 4746.2
$ws.Range("K40").Value = 4746.2
$ws.Range("M40").Value = -4610.2
$ws.Range("H46").Value = 3306.25
$ws.Range("I46").Value = 2992.8572
$ws.Range("K46").Value = 2992.8572
$ws.Range("M46").Value = -2804.8572
$ws.Range("H82").Value = 2416.3845
$ws.Range("I82").Value = 2429.2222
$ws.Range("J82").Value = 2387.5
$ws.Range("K82").Value = 2429.2222
$ws.Range("L82").Value = 2387.5
$ws.Range("M82").Value = -2068.2222
$ws.Range("N82").Value = -3109.5
$ws.Range("H85").Value = 2416.3845
$ws.Range("I85").Value = 2429.2222
$ws.Range("J85").Value = 2387.5
$ws.Range("K85").Value = 2429.2222
$ws.Range("L85").Value = 2387.5
$ws.Range("M85").Value = -1181.2222
$ws.Range("N85").Value = -4883.5
$ws.Range("I126").Value = 7007.0835
$ws.Range("J126").Value = 9499.5
$ws.Range("K126").Value = 21021.2505
$ws.Range("L126").Value = 28498.5
$ws.Range("M126").Value = -18551.2505
$ws.Range("N126").Value = -33438.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 83335230
$ws.Range("I122").Value = 2850
$ws.Range("J122").Value = 250000000
$ws.Range("K122").Value = 8550
$ws.Range("L122").Value = 750000000
$ws.Range("M122").Value = -6100
$ws.Range("N122").Value = -750004900
$ws.Range("H136").Value = 2057.7144
$ws.Range("I136").Value = 1878.8
$ws.Range("J136").Value = 2505
$ws.Range("K136").Value = 5636.4
$ws.Range("L136").Value = 7515
$ws.Range("M136").Value = -3086.4
$ws.Range("N136").Value = -12615
